# Update the dSF (column F) values for the pulled/repulled data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 2
    4  = 3
    5  = -2
    7  = 0
    8  = -3
    9  = 4
    10 = -1
    11 = -3
    13 = 2
    14 = 3
    15 = 7
    16 = -2
    17 = 3
    18 = -6
    19 = -5
    20 = 1
    23 = 1
    24 = 2
    25 = -6
    27 = -1
    28 = -3
    29 = 1
    30 = -3
    31 = -1
    32 = 6
    34 = -4
    36 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
